$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.862658023834229
$ws.Range("B1").Value = 3.087390422821045
$ws.Range("C1").Value = 2.159234523773193
$ws.Range("D1").Value = 1.944765686988831
$ws.Range("E1").Value = 1.790371894836426
